$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 56
$ws1.Range("F3").Value = 7257
$ws1.Range("F4").Value = 3504
$ws1.Range("F6").Value = 3833
$ws1.Range("F9").Value = 74
$ws1.Range("F11").Value = 143
$ws1.Range("F12").Value = 505
$ws1.Range("F14").Value = 130
$ws1.Range("F15").Value = 362
$ws1.Range("F19").Value = 4104
$ws1.Range("F21").Value = 408
$ws1.Range("F22").Value = 1027
$ws1.Range("F23").Value = 534
$ws1.Range("F24").Value = 1646
$ws1.Range("F27").Value = 3013
$ws1.Range("F28").Value = 2202
$ws1.Range("F32").Value = 22
$ws1.Range("F33").Value = 89
$ws1.Range("F36").Value = 4277
$ws1.Range("F37").Value = 467
$ws1.Range("F41").Value = 788
$ws1.Range("F42").Value = 200
$ws1.Range("F44").Value = 1624
$ws1.Range("F48").Value = 714

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 568

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 56
$ws4.Range("F5").Value = 7257
$ws4.Range("F6").Value = 3504
$ws4.Range("F7").Value = 3504
$ws4.Range("F8").Value = 3833
$ws4.Range("F10").Value = 74
$ws4.Range("F13").Value = 143
$ws4.Range("F14").Value = 505
$ws4.Range("F16").Value = 130
$ws4.Range("F17").Value = 362
$ws4.Range("F21").Value = 4104
$ws4.Range("F25").Value = 408
$ws4.Range("F26").Value = 1027
$ws4.Range("F27").Value = 534
$ws4.Range("F28").Value = 1646
$ws4.Range("F31").Value = 3013
$ws4.Range("F32").Value = 2202
$ws4.Range("F39").Value = 4277
$ws4.Range("F41").Value = 467
$ws4.Range("F50").Value = 714
# --- Sheet "all types": row 43 (old GOJO-15th event) content is replaced by what was
#     row 44 content; each of rows 44-47 takes on the next row old content (shift up
#     by one), and a brand-new event (EXA 2nd Diduhaihe-only) becomes row 48. ---
$ws4.Range("C43").Value = "北京·LookLook剧情式沉浸游戏互动动漫嘉年华"
$ws4.Range("D43").Value = "东村文化创意产业园A1-2 五道杠实景片场"
$ws4.Range("E43").Value = "2024.08.10 09:30-08.11 17:30"
$ws4.Range("F43").Value = 941
$ws4.Range("G43").Value = 29.9
$ws4.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=84741"
$ws4.Range("I43").Value = "//i2.hdslb.com/bfs/openplatform/202405/WH4KKudj1716880619473.jpeg"

$ws4.Range("C44").Value = "北京·梦次元动漫展M30"
$ws4.Range("D44").Value = "北京展览馆 北京展览馆"
$ws4.Range("E44").Value = "2024.08.10 10:00-08.11 17:00"
$ws4.Range("F44").Value = 788
$ws4.Range("G44").Value = 80
$ws4.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=83828"
$ws4.Range("I44").Value = "//i1.hdslb.com/bfs/openplatform/202405/Qr2Bd5W41715931423636.jpeg"

$ws4.Range("C45").Value = "北京·第五人格ONLY2.0"
$ws4.Range("D45").Value = "永外高庄138号 北京大红门国际会展中心"
$ws4.Range("E45").Value = "2024.08.10 10:00-08.10 17:00"
$ws4.Range("F45").Value = 200
$ws4.Range("G45").Value = 60
$ws4.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=86590"
$ws4.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202405/4jQBoo241716968548735.jpeg"

$ws4.Range("B46").NumberFormat = "@"
$ws4.Range("B46").Value = "2024-08-17"
$ws4.Range("B46").Style = "Normal"
$ws4.Range("C46").Value = "北京·第18届IJOY漫展xCGF游戏节"
$ws4.Range("D46").Value = "天辰东路7号 北京国家会议中心"
$ws4.Range("E46").Value = "2024.08.17 09:00-08.18 17:00"
$ws4.Range("F46").Value = 1624
$ws4.Range("G46").Value = 85
$ws4.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=84128"
$ws4.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202405/TU8kiduQ1715238040248.jpeg"

$ws4.Range("B47").NumberFormat = "@"
$ws4.Range("B47").Value = "2024-08-24"
$ws4.Range("B47").Style = "Normal"
$ws4.Range("C47").Value = "北京·万游引力夏日动漫游戏狂欢节"
$ws4.Range("D47").Value = "北七家镇王府街55号 水城会议中心"
$ws4.Range("E47").Value = "2024.08.24 10:00-08.24 17:00"
$ws4.Range("F47").Value = 260
$ws4.Range("G47").Value = 75
$ws4.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=83880"
$ws4.Range("I47").Value = "//i1.hdslb.com/bfs/openplatform/202404/vfXK7QCz1712541413788.jpeg"

$ws4.Range("B48").NumberFormat = "@"
$ws4.Range("B48").Value = "2024-09-15"
$ws4.Range("B48").Style = "Normal"
$ws4.Range("C48").Value = "北京· EXA·第二届帝都百合only"
$ws4.Range("D48").Value = "永外高庄138号 大红门国际会展中心"
$ws4.Range("E48").Value = "2024.09.15 09:30-09.15 16:00"
$ws4.Range("F48").Value = 29
$ws4.Range("G48").Value = 68
$ws4.Range("H48").Value = "https://show.bilibili.com/platform/detail.html?id=86477"
$ws4.Range("I48").Value = "//i2.hdslb.com/bfs/openplatform/202405/LgmTjud21716883789133.jpeg"
